# The "Material Status" column (M) and "Transport details" column (O) were
# separated by an unused, empty column L, and followed by another unused,
# empty column N. Remove those two empty spacer columns so the data packs
# together: M ("Material Status") shifts left into L, and O ("Transport
# details") shifts left into M.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("L").Delete()
$ws.Columns("M").Delete()
